$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) to remain text while we overwrite it, so plain
# decimal-looking values ("0.9964", "322.25", ...) are not silently
# reinterpreted as numbers by Excel (matches source inlineStr cells).
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "27.936.41"
$ws.Range("E2").Value = "  +1.15%  "
$ws.Range("D3").Value = "1.774.12"
$ws.Range("E3").Value = "  +0.73%  "
$ws.Range("D4").Value = "0.9964"
$ws.Range("E4").Value = "  -0.93%  "
$ws.Range("D5").Value = "322.25"
$ws.Range("E5").Value = "  -0.95%  "
$ws.Range("D6").Value = "0.9952"
$ws.Range("E6").Value = "  -0.73%  "
$ws.Range("D7").Value = "0.4272"
$ws.Range("E7").Value = "  -5.27%  "
$ws.Range("D8").Value = "0.3616"
$ws.Range("E8").Value = "  -2.76%  "
$ws.Range("D9").Value = "44.41"
$ws.Range("E9").Value = "  -1.91%  "
$ws.Range("E10").Value = "  -3.61%  "
$ws.Range("D11").Value = "1.107"
$ws.Range("E11").Value = "  -1.62%  "
$ws.Range("E12").Value = "  -0.90%  "
$ws.Range("D13").Value = "21.59"
$ws.Range("E13").Value = "  -0.81%  "
$ws.Range("D14").Value = "6.131"
$ws.Range("E14").Value = "  -0.97%  "
$ws.Range("D15").Value = "7.320"
$ws.Range("E15").Value = "  -0.83%  "
$ws.Range("D16").Value = "1.788.54"
$ws.Range("E16").Value = "  +1.52%  "
$ws.Range("D17").Value = "91.50"
$ws.Range("E17").Value = "  -0.09%  "
$ws.Range("D18").Value = "0.00001061"
$ws.Range("E18").Value = "  -1.80%  "
$ws.Range("D19").Value = "0.06342"
$ws.Range("E19").Value = "  +0.82%  "
$ws.Range("D20").Value = "0.9952"
$ws.Range("E20").Value = "  -0.66%  "
$ws.Range("D21").Value = "17.20"
$ws.Range("E21").Value = "  -1.36%  "
$ws.Range("D22").Value = "5.950"
$ws.Range("E22").Value = "  -3.75%  "
$ws.Range("D23").Value = "27.935.07"
$ws.Range("E23").Value = "  +0.96%  "
$ws.Range("E24").Value = "  -2.16%  "
$ws.Range("D25").Value = "2.157"
$ws.Range("E25").Value = "  -6.58%  "
$ws.Range("D26").Value = "159.98"
$ws.Range("E26").Value = "  +4.13%  "
$ws.Range("E27").Value = "  -2.25%  "
$ws.Range("D28").Value = "1.989.15"
$ws.Range("E28").Value = "  +1.13%  "
$ws.Range("D29").Value = "2.167"
$ws.Range("E29").Value = "  -6.92%  "
$ws.Range("D30").Value = "126.09"
$ws.Range("E30").Value = "  -2.16%  "
$ws.Range("D31").Value = "1.167"
$ws.Range("E31").Value = "  -2.94%  "
$ws.Range("D32").Value = "5.688"
$ws.Range("E32").Value = "  -1.43%  "
$ws.Range("D33").Value = "0.08977"
$ws.Range("E33").Value = "  -2.70%  "
$ws.Range("D34").Value = "3.503"
$ws.Range("E34").Value = "  -4.83%  "
$ws.Range("D35").Value = "12.64"
$ws.Range("E35").Value = "  -1.21%  "
$ws.Range("E36").Value = "  -0.46%  "
$ws.Range("D37").Value = "5.067"
$ws.Range("E37").Value = "  -0.29%  "
$ws.Range("D38").Value = "0.2113"
$ws.Range("E38").Value = "  -3.15%  "
$ws.Range("D39").Value = "0.6433"
$ws.Range("E39").Value = "  -0.64%  "
$ws.Range("D40").Value = "0.06051"
$ws.Range("E40").Value = "  -1.07%  "
$ws.Range("E41").Value = "  -0.83%  "
$ws.Range("D42").Value = "0.9946"
$ws.Range("E42").Value = "  -0.72%  "
$ws.Range("D43").Value = "7.855"
$ws.Range("E43").Value = "  -2.12%  "
$ws.Range("D44").Value = "1.389"
$ws.Range("E44").Value = "  -1.31%  "
$ws.Range("E45").Value = "  -1.64%  "
$ws.Range("D46").Value = "0.5970"
$ws.Range("E46").Value = "  -0.49%  "
$ws.Range("D47").Value = "3.689"
$ws.Range("E47").Value = "  -1.46%  "
$ws.Range("D48").Value = "124.28"
$ws.Range("E48").Value = "  -1.28%  "
$ws.Range("D49").Value = "1.983"
$ws.Range("E49").Value = "  -0.60%  "
$ws.Range("D50").Value = "1.150"
$ws.Range("E50").Value = "  +0.42%  "
$ws.Range("D51").Value = "0.06898"
$ws.Range("E51").Value = "  -0.22%  "

# Restore the default (unstyled) cell format so the written cells match
# the original workbook's style (no explicit style index).
$ws.Range("D2:D51").Style = "Normal"

